# Insert a new weekly price-report row above the current first data row (row 9).
# This pushes the existing rows 9:99 down to 10:100 (dimension grows to A1:R100)
# and the new row 9 gets a fresh set of observations for this "Hortaliza /
# Berenjena - Vega Monumental Concepción" consolidation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("9:9").Insert()

$ws.Range("A9").Value = 11
$ws.Range("B9").Value = "Vega Monumental Concepción"
$ws.Range("C9").Value = "Bíobío"
$ws.Range("D9").Value = 44750
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 100112001
$ws.Range("G9").Value = "Berenjena"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 110
$ws.Range("K9").Value = 9500
$ws.Range("L9").Value = 10000
$ws.Range("M9").Value = 9727
$ws.Range("N9").Value = "`$/caja 60 unidades"
$ws.Range("O9").Value = "Región de Arica y Parinacota"
$ws.Range("P9").Value = 162
$ws.Range("Q9").Value = 60
$ws.Range("R9").Value = "Hortaliza"
